$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92

# Text columns (A-D) would otherwise be auto-converted by Excel's type
# inference (dates, leading-zero numbers, etc.), so force text storage
# via a temporary "@" number format, then clear formatting so the cell
# ends up styled like its neighbours (no explicit style index).
$textCells = @{
    1 = "2024-01-24"
    2 = "22:08:36"
    3 = "Wednesday"
    4 = "03"
}
foreach ($col in $textCells.Keys) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$col]
    $cell.ClearFormats()
}

# Numeric columns (E-T)
$numCells = @{
    5  = 138559
    6  = 141498
    7  = 171515
    8  = 149058
    9  = -1
    10 = 123515
    11 = 223888
    12 = 256370
    13 = 185149
    14 = 110045
    15 = 41320
    16 = 30892
    17 = 73464
    18 = -1
    19 = 42501
    20 = -1
}
foreach ($col in $numCells.Keys) {
    $ws.Cells.Item($row, $col).Value = $numCells[$col]
}
